$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.924.08"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "2.636.04"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.45"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.57"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").Value = "2.678.27"
$ws.Range("E9").Value = "  +1.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.52"
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.339"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "3.105.41"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "58.870.83"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.29"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "2.667.92"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.60"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.15"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.51"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.91"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.425"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "2.765.94"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").Value = "0.0₃0825"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.10"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.59"
$ws.Range("E32").Value = "  +8.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.01"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.58"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.09"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("E36").Value = "  +15.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.03"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.15"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.864"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.56"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.69"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.42"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.630"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "281.65"
$ws.Range("E44").Value = "  -2.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0994"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0542"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.78"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0232"
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.24"
$ws.Range("E51").Value = "  -1.04%  "
